# "drop down stat-list added" -- append a new block of answers (unique_key
# 99295) to the "Ответы" sheet, one row per standard question, all pointing
# at the PVI location "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ответы")

# Seed rows 11-13 by copying the formatting of the last existing data row
# (row 10) so the new cells pick up the same border/bold/centered style for
# column A and the same date number format for column F, without minting
# any new style entries.
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G13").PasteSpecial(-4122)

$questions = @(
    "Довольны ли Вы качеством предоставляемого питания?",
    "Устраивают ли Вас бытовые условия? (питьевой режим, температура в помещении, досуг и психологический климат)",
    "Довольны ли Вы работой обслуживающего персонала?"
)

$uniqueKey = 99295
$pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)"

$startRow = 11
$idx = 9
for ($i = 0; $i -lt $questions.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $uniqueKey
    $ws.Cells.Item($r, 3).Value = $questions[$i]
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = "None"
    $ws.Cells.Item($r, 6).Value = 44165
    $ws.Cells.Item($r, 7).Value = $pvi
    $idx = $idx + 1
}
